$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33, shifting existing rows 33:121 down to 34:122
$ws.Rows("33").Insert()

# Populate the newly inserted row 33 with the new record (same static
# fields as the template rows, new Fecha/Volumen/Precio values)
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44622
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = 100112030
$ws.Range("G33").Value = "Poroto granado"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 23000
$ws.Range("L33").Value = 23000
$ws.Range("M33").Value = 23000
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 920
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
